# Apply updates to the quarterly dollar income statement workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the period-end date labels (row 9) reflecting the new report date.
$ws.Range("I9").Value = "1402-03-13 (10)"
$ws.Range("M9").Value = "1402-03-13 (2)"

# Update the latest quarter (column M) figures.
$ws.Range("M14").Value = -7576
$ws.Range("M17").Value = 6627
$ws.Range("M18").Value = -5468
$ws.Range("M20").Value = 46018
$ws.Range("M21").Value = 305
$ws.Range("M22").Value = 46323
$ws.Range("M24").Value = 46323
